$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "290.56"
Set-TextValue $ws.Range("E2") "-3.70%"

Set-TextValue $ws.Range("D3") "30.80"
Set-TextValue $ws.Range("E3") "-6.11%"

Set-TextValue $ws.Range("D4") "4.959"
Set-TextValue $ws.Range("E4") "0.20%"

Set-TextValue $ws.Range("E5") "-7.07%"

Set-TextValue $ws.Range("D6") "1.792"
Set-TextValue $ws.Range("E6") "-8.80%"

Set-TextValue $ws.Range("D7") "7.680"
Set-TextValue $ws.Range("E7") "-2.19%"

Set-TextValue $ws.Range("D8") "3.764"
Set-TextValue $ws.Range("E8") "-0.94%"

Set-TextValue $ws.Range("D9") "0.8960"
Set-TextValue $ws.Range("E9") "-3.06%"

Set-TextValue $ws.Range("D10") "0.1654"
Set-TextValue $ws.Range("E10") "-6.25%"

Set-TextValue $ws.Range("D11") "0.07728"
Set-TextValue $ws.Range("E11") "-1.69%"

Set-TextValue $ws.Range("D12") "0.08068"
Set-TextValue $ws.Range("E12") "-6.21%"

Set-TextValue $ws.Range("D13") "0.03041"
Set-TextValue $ws.Range("E13") "-3.44%"

Set-TextValue $ws.Range("E14") "0.16%"

Set-TextValue $ws.Range("D15") "0.001504"
Set-TextValue $ws.Range("E15") "-0.57%"

Set-TextValue $ws.Range("D16") "0.005756"
Set-TextValue $ws.Range("E16") "-2.42%"

Set-TextValue $ws.Range("D17") "3.468"
Set-TextValue $ws.Range("E17") "0.17%"

Set-TextValue $ws.Range("D18") "2.082"
Set-TextValue $ws.Range("E18") "-3.33%"

Set-TextValue $ws.Range("D19") "0.3312"
Set-TextValue $ws.Range("E19") "-0.74%"

Set-TextValue $ws.Range("D20") "0.1302"
Set-TextValue $ws.Range("E20") "-1.25%"

Set-TextValue $ws.Range("E21") "-5.99%"

Set-TextValue $ws.Range("D22") "0.2323"
Set-TextValue $ws.Range("E22") "16.57%"

Set-TextValue $ws.Range("D23") "0.04510"
Set-TextValue $ws.Range("E23") "-0.84%"

Set-TextValue $ws.Range("E24") "-0.97%"

Set-TextValue $ws.Range("D25") "0.004014"
Set-TextValue $ws.Range("E25") "-9.42%"

Set-TextValue $ws.Range("D26") "0.0001252"
Set-TextValue $ws.Range("E26") "0.02%"

Set-TextValue $ws.Range("D39") "0.01604"
Set-TextValue $ws.Range("E39") "-6.34%"

Set-TextValue $ws.Range("D40") "0.04396"
Set-TextValue $ws.Range("E40") "-6.99%"

Set-TextValue $ws.Range("D41") "0.007287"
Set-TextValue $ws.Range("E41") "-9.99%"

Set-TextValue $ws.Range("E42") "-3.30%"

Set-TextValue $ws.Range("D43") "0.007682"

Set-TextValue $ws.Range("D44") "0.001902"
Set-TextValue $ws.Range("E44") "-18.79%"

Set-TextValue $ws.Range("D45") "0.009213"
Set-TextValue $ws.Range("E45") "-12.56%"

Set-TextValue $ws.Range("D46") "0.00005941"
Set-TextValue $ws.Range("E46") "-5.17%"

Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("E47") "0.01%"

Set-TextValue $ws.Range("D48") "2.246"
Set-TextValue $ws.Range("E48") "172.73%"

Set-TextValue $ws.Range("E49") "-3.21%"

Set-TextValue $ws.Range("D50") "0.00002102"
Set-TextValue $ws.Range("E50") "0.01%"

Set-TextValue $ws.Range("D51") "0.0002002"
Set-TextValue $ws.Range("E51") "0.01%"

Write-Host "Applied all cryptos price/volume updates"